# M10 Data Aug Froze E1
# Updates the epoch-accuracy figures in column B (re-run results) and
# refreshes the repr() memory-address text stored in column A for the
# DisplayOutputs rows (102-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: updated per-epoch accuracy values -------------------------
$ws.Cells.Item(2, 2).Value = 0.265625
$ws.Cells.Item(3, 2).Value = 0.140625
$ws.Cells.Item(4, 2).Value = 0.125
$ws.Cells.Item(5, 2).Value = 0.140625
$ws.Cells.Item(6, 2).Value = 0.15625
$ws.Cells.Item(7, 2).Value = 0.1875
$ws.Cells.Item(8, 2).Value = 0.15625
$ws.Cells.Item(9, 2).Value = 0.140625
$ws.Cells.Item(10, 2).Value = 0.125
$ws.Cells.Item(13, 2).Value = 0.140625
$ws.Cells.Item(14, 2).Value = 0.3125
$ws.Cells.Item(16, 2).Value = 0.15625
$ws.Cells.Item(17, 2).Value = 0.109375
$ws.Cells.Item(18, 2).Value = 0.140625
$ws.Cells.Item(19, 2).Value = 0.171875
$ws.Cells.Item(20, 2).Value = 0.40625
$ws.Cells.Item(21, 2).Value = 0.203125
$ws.Cells.Item(22, 2).Value = 0.109375
$ws.Cells.Item(23, 2).Value = 0.21875
$ws.Cells.Item(24, 2).Value = 0.078125
$ws.Cells.Item(25, 2).Value = 0.125
$ws.Cells.Item(26, 2).Value = 0.125
$ws.Cells.Item(27, 2).Value = 0.109375
$ws.Cells.Item(29, 2).Value = 0.109375
$ws.Cells.Item(30, 2).Value = 0.109375
$ws.Cells.Item(31, 2).Value = 0.109375
$ws.Cells.Item(32, 2).Value = 0.109375
$ws.Cells.Item(33, 2).Value = 0.109375
$ws.Cells.Item(34, 2).Value = 0.109375
$ws.Cells.Item(35, 2).Value = 0.109375
$ws.Cells.Item(36, 2).Value = 0.109375
$ws.Cells.Item(37, 2).Value = 0.09375
$ws.Cells.Item(38, 2).Value = 0.09375
$ws.Cells.Item(39, 2).Value = 0.09375
$ws.Cells.Item(40, 2).Value = 0.09375
$ws.Cells.Item(41, 2).Value = 0.09375
$ws.Cells.Item(42, 2).Value = 0.09375
$ws.Cells.Item(43, 2).Value = 0.09375
$ws.Cells.Item(44, 2).Value = 0.09375
$ws.Cells.Item(45, 2).Value = 0.09375
$ws.Cells.Item(46, 2).Value = 0.09375
$ws.Cells.Item(47, 2).Value = 0.09375
$ws.Cells.Item(48, 2).Value = 0.09375
$ws.Cells.Item(49, 2).Value = 0.09375
$ws.Cells.Item(50, 2).Value = 0.09375
$ws.Cells.Item(51, 2).Value = 0.09375
$ws.Cells.Item(52, 2).Value = 0.09375
$ws.Cells.Item(53, 2).Value = 0.09375
$ws.Cells.Item(54, 2).Value = 0.09375
$ws.Cells.Item(55, 2).Value = 0.09375
$ws.Cells.Item(56, 2).Value = 0.09375
$ws.Cells.Item(57, 2).Value = 0.09375
$ws.Cells.Item(58, 2).Value = 0.09375
$ws.Cells.Item(59, 2).Value = 0.09375
$ws.Cells.Item(60, 2).Value = 0.09375
$ws.Cells.Item(61, 2).Value = 0.09375
$ws.Cells.Item(62, 2).Value = 0.09375
$ws.Cells.Item(63, 2).Value = 0.109375
$ws.Cells.Item(64, 2).Value = 0.109375
$ws.Cells.Item(65, 2).Value = 0.109375
$ws.Cells.Item(103, 2).Value = 0.125
$ws.Cells.Item(105, 2).Value = 0.0625
$ws.Cells.Item(106, 2).Value = 0.0625
$ws.Cells.Item(109, 2).Value = 0.078125
$ws.Cells.Item(110, 2).Value = 0.0625
$ws.Cells.Item(111, 2).Value = 0.03125
$ws.Cells.Item(112, 2).Value = 0.046875
$ws.Cells.Item(113, 2).Value = 0.15625
$ws.Cells.Item(114, 2).Value = 0.078125
$ws.Cells.Item(117, 2).Value = 0.0625
$ws.Cells.Item(118, 2).Value = 0.1475409836065574

# --- Column A (rows 102-118): refresh the DisplayOutputs repr() text -----
# The object's id()/memory address changed between the old run and the new
# run, so every occurrence of the old hex address is replaced with the new
# one across all affected rows.
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = "<__main__.DisplayOutputs object at 0x7f79dfef0fd0>"
}
